$wb = $excel.ActiveWorkbook

# ---------- GLOBAL RESULTS ----------
$ws = $wb.Worksheets.Item("GLOBAL RESULTS")
$ws.Range("C6").Value  = 51452.94134312038
$ws.Range("C7").Value  = 51452.94134312038
$ws.Range("C8").Value  = 46307.64720880834
$ws.Range("C12").Value = 39493.87466383973
$ws.Range("C13").Value = 39493.87466383973
$ws.Range("C14").Value = 26623.874663839728
$ws.Range("C15").Value = 25894.787372839724
$ws.Range("C16").Value = 25044.577372839733
$ws.Range("C20").Value = 504580.9872225114
$ws.Range("C21").Value = 504580.9872225114
$ws.Range("C22").Value = 454122.8885002602
$ws.Range("C26").Value = 387302.60597214376
$ws.Range("C27").Value = 387302.60597214376
$ws.Range("C28").Value = 261091.02047214383
$ws.Range("C29").Value = 253941.11658985866
$ws.Range("C30").Value = 245603.40469335875

# ---------- FUSELAGE ----------
$ws = $wb.Worksheets.Item("FUSELAGE")
$ws.Range("C7").Value  = 5710.0
$ws.Range("D7").Value  = 14.09275281235641
$ws.Range("C8").Value  = 5671.0
$ws.Range("D8").Value  = 13.313485323795659
$ws.Range("C9").Value  = 5985.0
$ws.Range("D9").Value  = 19.58758766759249
$ws.Range("C12").Value = 6402.333333333332
$ws.Range("D12").Value = 27.926415835781025

# ---------- WING ----------
$ws = $wb.Worksheets.Item("WING")
$ws.Range("C7").Value  = 5078.0
$ws.Range("D7").Value  = 58.23751207503664
$ws.Range("C8").Value  = 4036.0
$ws.Range("D8").Value  = 25.76734910099407
$ws.Range("C9").Value  = 5097.0
$ws.Range("D9").Value  = 58.829578386463524
$ws.Range("C11").Value = 5224.0
$ws.Range("D11").Value = 62.78707425758004
$ws.Range("C12").Value = 4676.0
$ws.Range("D12").Value = 45.71063538063634
$ws.Range("C13").Value = 4188.428571428571
$ws.Range("D13").Value = 30.517234471614202

# ---------- HORIZONTAL TAIL ----------
$ws = $wb.Worksheets.Item("HORIZONTAL TAIL")
$ws.Range("C8").Value  = 235.0
$ws.Range("D8").Value  = -57.140251687032645
$ws.Range("C9").Value  = 440.0
$ws.Range("D9").Value  = -19.75196060550793
$ws.Range("C10").Value = 468.0
$ws.Range("D10").Value = -14.645267189494787

# ---------- VERTICAL TAIL ----------
$ws = $wb.Worksheets.Item("VERTICAL TAIL")
$ws.Range("C8").Value = 404.0
$ws.Range("D8").Value = -10.757676165230823
$ws.Range("C9").Value = 471.5
$ws.Range("D9").Value = 4.15286061409322

# ---------- POWER PLANT ----------
$ws = $wb.Worksheets.Item("POWER PLANT")
$ws.Range("C2").Value  = 2699.999999999999
$ws.Range("C3").Value  = 4049.999999999999
$ws.Range("C8").Value  = 1350.0
$ws.Range("C9").Value  = 2024.9999999999995
$ws.Range("C12").Value = 1350.0
$ws.Range("C13").Value = 2024.9999999999995

# ---------- LANDING GEARS ----------
$ws = $wb.Worksheets.Item("LANDING GEARS")
$ws.Range("C5").Value = 1499.0
$ws.Range("D5").Value = -32.87358380726343
$ws.Range("C6").Value = 2055.0
$ws.Range("D6").Value = -7.975460122699364
$ws.Range("C7").Value = 2336.0
$ws.Range("D7").Value = 4.6079441135641295
$ws.Range("C8").Value = 2024.0
$ws.Range("D8").Value = -9.363664860507791
$ws.Range("C9").Value = 1978.5
$ws.Range("D9").Value = -11.401191169226626
